$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "SCD0278" to "SCD0018"
$ws.Name = "SCD0018"

# Update TC_ID column (B2:B4) value from "DGS-293" to "SCD0018-001"
$ws.Range("B2:B4").Value = "SCD0018-001"

# Column B needs to widen (best-fit) to accommodate the longer TC_ID text
$ws.Columns("B").ColumnWidth = (70.0/6.0)

# Move the active selection to B5 and scroll the view back to show column A
$ws.Range("B5").Select()
